$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 0.1363636363636364
$ws.Range("E2").Value2 = 0.6893090205863619
$ws.Range("D3").Value2 = 0.2676962676962677
$ws.Range("E3").Value2 = 0.1144610211151944
$ws.Range("D4").Value2 = 0.4209486166007905
$ws.Range("E4").Value2 = 0.04547120480280356
$ws.Range("D5").Value2 = -0.5000000000000001
$ws.Range("E5").Value2 = 0.20703125
$ws.Range("D6").Value2 = 0.7714285714285715
$ws.Range("E6").Value2 = 0.07239650145772594
$ws.Range("D7").Value2 = 0.4509803921568628
$ws.Range("E7").Value2 = 0.06032623578568472
$ws.Range("D8").Value2 = 0.1818181818181818
$ws.Range("E8").Value2 = 0.5926152128455
$ws.Range("D9").Value2 = 0.2272844272844273
$ws.Range("E9").Value2 = 0.1824978117098565
$ws.Range("D10").Value2 = 0.3893280632411067
$ws.Range("E10").Value2 = 0.06632062116006986
$ws.Range("D11").Value2 = -0.5000000000000001
$ws.Range("E11").Value2 = 0.20703125
$ws.Range("D12").Value2 = 0.7714285714285715
$ws.Range("E12").Value2 = 0.07239650145772594
$ws.Range("D13").Value2 = 0.4509803921568628
$ws.Range("E13").Value2 = 0.06032623578568472
$ws.Range("D14").Value2 = 0.2181818181818182
$ws.Range("E14").Value2 = 0.519248247554982
$ws.Range("D15").Value2 = 0.1552123552123552
$ws.Range("E15").Value2 = 0.3660469538335063
$ws.Range("D16").Value2 = 0.4189723320158103
$ws.Range("E16").Value2 = 0.04660153881127199
$ws.Range("D17").Value2 = -0.5000000000000001
$ws.Range("E17").Value2 = 0.20703125
$ws.Range("D18").Value2 = 0.7714285714285715
$ws.Range("E18").Value2 = 0.07239650145772594
$ws.Range("D19").Value2 = 0.4509803921568628
$ws.Range("E19").Value2 = 0.06032623578568472
$ws.Range("D20").Value2 = 0.2272727272727273
$ws.Range("E20").Value2 = 0.501535667604937
$ws.Range("D21").Value2 = 0.1158301158301158
$ws.Range("E21").Value2 = 0.5011231999211077
$ws.Range("D22").Value2 = 0.3784584980237154
$ws.Range("E22").Value2 = 0.07495820299718052
$ws.Range("D23").Value2 = -0.5000000000000001
$ws.Range("E23").Value2 = 0.20703125
$ws.Range("D24").Value2 = 0.7714285714285715
$ws.Range("E24").Value2 = 0.07239650145772594
$ws.Range("D25").Value2 = 0.4509803921568628
$ws.Range("E25").Value2 = 0.06032623578568472
$ws.Range("D26").Value2 = 0.2090909090909091
$ws.Range("E26").Value2 = 0.5372209352113229
$ws.Range("D27").Value2 = 0.07413127413127413
$ws.Range("E27").Value2 = 0.6674266915157008
$ws.Range("D28").Value2 = 0.3853754940711462
$ws.Range("E28").Value2 = 0.06936888423715426
$ws.Range("D29").Value2 = -0.5000000000000001
$ws.Range("E29").Value2 = 0.20703125
$ws.Range("D30").Value2 = 0.7714285714285715
$ws.Range("E30").Value2 = 0.07239650145772594
$ws.Range("D31").Value2 = 0.4509803921568628
$ws.Range("E31").Value2 = 0.06032623578568472
$ws.Range("D32").Value2 = 0.3090909090909091
$ws.Range("E32").Value2 = 0.3550284397015071
$ws.Range("D33").Value2 = 0.0332046332046332
$ws.Range("E33").Value2 = 0.8475467370962678
$ws.Range("D34").Value2 = 0.3567193675889328
$ws.Range("E34").Value2 = 0.09475100941564729
$ws.Range("D35").Value2 = -0.5000000000000001
$ws.Range("E35").Value2 = 0.20703125
$ws.Range("D36").Value2 = 0.7714285714285715
$ws.Range("E36").Value2 = 0.07239650145772594
$ws.Range("D37").Value2 = 0.4509803921568628
$ws.Range("E37").Value2 = 0.06032623578568472

Write-Output "Updated correlation and p-value columns for rows 2-37"
